# Update ligand/receptor TPM-derived expression values in the NATMI
# LR-pairs output sheet (Plau -> Igf2r), reflecting the new TPM values
# used upstream ("update scripts wuth new tpm").
#
# New ligand (Plau) average/total expression values per Sending cluster
# and new receptor (Igf2r) average/total expression values per Target
# cluster. All other dependent columns (specificity, edge weights, edge
# specificity) are recomputed from these.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Ligand average / total expression value, keyed by Sending cluster
$sendingLigand = @{
    "ECs"              = @(27.67568533333333, 83.027056)
    "FAPs"             = @(44.54713066666667, 133.641392)
    "Inflammatory-Mac" = @(54.059897, 162.179691)
    "MuSCs"            = @(11.41370466666667, 34.241114)
    "Resolving-Mac"    = @(44.88488133333333, 134.654644)
}

# New Receptor average / total expression value, keyed by Target cluster
$targetReceptor = @{
    "ECs"              = @(30.63045066666666, 91.891352)
    "FAPs"             = @(52.681366, 158.044098)
    "Inflammatory-Mac" = @(6.676334999999999, 20.029005)
    "MuSCs"            = @(20.21956533333333, 60.658696)
    "Resolving-Mac"    = @(11.35077233333333, 34.052317)
}

# Sum across all clusters, used to derive the specificity columns
$totalG = 0.0
foreach ($key in $sendingLigand.Keys) { $totalG += $sendingLigand[$key][0] }
$totalH = 0.0
foreach ($key in $sendingLigand.Keys) { $totalH += $sendingLigand[$key][1] }
$totalM = 0.0
foreach ($key in $targetReceptor.Keys) { $totalM += $targetReceptor[$key][0] }
$totalN = 0.0
foreach ($key in $targetReceptor.Keys) { $totalN += $targetReceptor[$key][1] }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 26 }

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target = $ws.Cells.Item($r, 4).Value2

    if (-not $sendingLigand.ContainsKey($sending)) { continue }
    if (-not $targetReceptor.ContainsKey($target)) { continue }

    $G = $sendingLigand[$sending][0]
    $H = $sendingLigand[$sending][1]
    $M = $targetReceptor[$target][0]
    $N = $targetReceptor[$target][1]

    $I = $G / $totalG
    $J = $H / $totalH
    $O = $M / $totalM
    $P = $N / $totalN

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J

    $ws.Cells.Item($r, 13).Value = $M
    $ws.Cells.Item($r, 14).Value = $N
    $ws.Cells.Item($r, 15).Value = $O
    $ws.Cells.Item($r, 16).Value = $P

    $ws.Cells.Item($r, 17).Value = $Q
    $ws.Cells.Item($r, 18).Value = $R
    $ws.Cells.Item($r, 19).Value = $S
    $ws.Cells.Item($r, 20).Value = $T
}
